$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Top providers" table written to columns Q:R, rows 2-12, mirroring the
# existing F:G table's layout/formatting (label left-aligned, count plain).

$data = @(
    @{ Row = 2;  Name = "YASIR  FASIH  MD";       Count = 7 },
    @{ Row = 3;  Name = "MARK A FELDNER  ";       Count = 5 },
    @{ Row = 4;  Name = "AMBIKA  RAO  ";          Count = 5 },
    @{ Row = 5;  Name = "IBRAHIM G ZABANEH  ";    Count = 5 },
    @{ Row = 6;  Name = "JOHN A HOEHN  ";         Count = 5 },
    @{ Row = 7;  Name = "CONSTANTINA C TUCKER  "; Count = 4 },
    @{ Row = 8;  Name = "AGNES  ERRIHANI  ";      Count = 3 },
    @{ Row = 9;  Name = "MANJEET  GEETA  MD";     Count = 3 },
    @{ Row = 10; Name = "ANAND A SHAH  ";         Count = 3 },
    @{ Row = 11; Name = "STEPHANIE A BRYANT  ";   Count = 3 },
    @{ Row = 12; Name = "JULIE A MORA  ";         Count = 3 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("Q$r").Value = $entry.Name
    $ws.Range("Q$r").HorizontalAlignment = -4131
    $ws.Range("R$r").Value = $entry.Count
}

$ws.Range("Q2:R12").Select()
